# Apply the "Contribution Summary" update:
#  - Replace the 4 placeholder task rows with the new, expanded contribution
#    log (14 rows covering Gabriel Tharp, Owen Randolph and Marcos Fernandez).
#  - Widen columns A and B to fit the new, longer names.
#  - Restore the cursor/selection and page orientation seen in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name, Task, Contribution Details, Hours
$data = @(
    @("Gabriel Tharp",    "Create MySQL DB",             "Create db and empty table in MySQL workbench with appropriate data types, constraints - SQL", 0.5),
    @("Gabriel Tharp",    "Clean/Transform",              "Remove duplicates, handle NULLs, clean data - Python", 1.5),
    @("Gabriel Tharp",    "Load Data",                    "Load CSV into Jupyter Notebook - Python", 0.1),
    @("Gabriel Tharp",    "ER Diagram",                   "Reverse engineer ER diagram - MySQL Workbench", 0.1),
    @("Gabriel Tharp",    "Load Data",                    "Make a connection into MySQL Workbench - Python", 0.25),
    @("Owen Randolph",    "Reporting",                    "Create structured draft of report", 3),
    @("Owen Randolph",    "Normalization",                "Check for 1NF Normalization - Python via ChatGPT", 0.5),
    @("Owen Randolph",    "Normalization",                "Create 2NF tables for normalization - SQL script", 1),
    @("Owen Randolph",    "Load Data",                    "Add data to new tables - SQL script", 0.25),
    @("Owen Randolph",    "Normalization",                "Add foreign keys for 3NF - SQL Script", 0.5),
    @("Marcos Fernandez", "Load Data",                    "Loaded data from jupyter Notebook to MySQL Workbench", 0.25),
    @("Marcos Fernandez", "Check for dulplicate values",  " After data transfer from csv in Jupyter directory to MySQL workbench - Python", 0.25),
    @("Marcos Fernandez", "Normalization",                "Create 3NF tables for normalization - SQL script", 1),
    @("Marcos Fernandez", "Reporting",                    "Refine Report Draft ", 3)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Columns A and B need to grow to fit the longer names / task labels.
$ws.Columns.Item(1).ColumnWidth = 16.67
$ws.Columns.Item(2).ColumnWidth = 24.17

# Restore the saved selection/cursor position.
[void]$ws.Range("B21").Select()

# Page was switched to (explicit) portrait orientation.
$ws.PageSetup.Orientation = 1
